$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '69.076.26'
$ws.Range('E2').Value = '  +1.97%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.817.21'
$ws.Range('E3').Value = '  +0.60%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.17%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '629.82'
$ws.Range('E5').Value = '  +5.24%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '165.40'
$ws.Range('E6').Value = '  +0.17%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.815.12'
$ws.Range('E7').Value = '  +0.63%  '
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.520'
$ws.Range('E9').Value = '  +0.64%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.161'
$ws.Range('E10').Value = '  +1.82%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.455'
$ws.Range('E11').Value = '  +0.89%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '6.60'
$ws.Range('E12').Value = '  +2.65%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000251'
$ws.Range('E13').Value = '  +1.02%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '36.05'
$ws.Range('E14').Value = '  +0.88%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.458.79'
$ws.Range('E15').Value = '  +0.61%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '4.028.58'
$ws.Range('E16').Value = '  +5.21%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '69.076.62'
$ws.Range('E17').Value = '  +1.91%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '18.05'
$ws.Range('E18').Value = '  -1.28%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.13'
$ws.Range('E19').Value = '  +1.20%  '
$ws.Range('E20').Value = '  -0.23%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '465.44'
$ws.Range('E21').Value = '  +0.93%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.68'
$ws.Range('E22').Value = '  -1.06%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.707'
$ws.Range('E23').Value = '  +1.29%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.0000152'
$ws.Range('E24').Value = '  +5.07%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '83.89'
$ws.Range('E25').Value = '  +1.37%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.96'
$ws.Range('E26').Value = '  -0.47%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.14'
$ws.Range('E27').Value = '  +2.24%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.02'
$ws.Range('E28').Value = '  +0.18%  '
$ws.Range('E29').Value = '  +0.01%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.971.28'
$ws.Range('E30').Value = '  +0.71%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.70'
$ws.Range('E31').Value = '  +1.98%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.22'
$ws.Range('E32').Value = '  +0.98%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '7.29'
$ws.Range('E33').Value = '  -1.75%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '29.13'
$ws.Range('E34').Value = '  +0.09%  '
$ws.Range('B35').Value = 'Binance-PegBSC-USD'
$ws.Range('C35').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.00'
$ws.Range('E35').Value = '  +0.11%  '
$ws.Range('B36').Value = 'Aptos'
$ws.Range('C36').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '9.08'
$ws.Range('E36').Value = '  +0.80%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.101'
$ws.Range('E37').Value = '  +2.29%  '
$ws.Range('E38').Value = '  +7.48%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.38'
$ws.Range('E39').Value = '  +5.26%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.91'
$ws.Range('E40').Value = '  +2.85%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.981'
$ws.Range('E41').Value = '  -0.82%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.00'
$ws.Range('E42').Value = '  +0.01%  '
$ws.Range('E43').Value = '  +0.03%  '
$ws.Range('B44').Value = 'Monero'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '156.80'
$ws.Range('E44').Value = '  +3.48%  '
$ws.Range('B45').Value = 'ONDO'
$ws.Range('C45').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.44'
$ws.Range('E45').Value = '  +4.97%  '
$ws.Range('B46').Value = 'TheGraph'
$ws.Range('C46').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.299'
$ws.Range('E46').Value = '  +0.47%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '46.73'
$ws.Range('E47').Value = '  -1.76%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.45'
$ws.Range('E48').Value = '  +1.34%  '
$ws.Range('B49').Value = 'Stacks'
$ws.Range('C49').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.89'
$ws.Range('E49').Value = '  +2.07%  '
$ws.Range('B50').Value = 'Arweave'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '42.16'
$ws.Range('E50').Value = '  -5.41%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.000280'
$ws.Range('E51').Value = '  +13.96%  '
